$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph: "03. Regras de Comunicação" -> "Regras de Comunicação"
#    (drop the "03. " numbering prefix)
# ------------------------------------------------------------------
$d.Content.Find.Execute("03. Regras de Comunicação", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Regras de Comunicação", 2)

# The old title also carried a leftover "_GoBack" bookmark spanning the
# "Regras de Com" / "unicação" run split - remove it now that the runs
# collapsed back into a single run of text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Body paragraph: drop "do grupo, " and the comma after "Davi"
#    "O representante do grupo, Davi, fará o contato..."
#      -> "O representante Davi fará o contato..."
# ------------------------------------------------------------------
$d.Content.Find.Execute("O representante do grupo, Davi, fará o contato", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "O representante Davi fará o contato", 2)

# ------------------------------------------------------------------
# 3) Body paragraph: trim the trailing sentence about feedback, keep
#    the trailing space that followed "pessoais."
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "reuniões pessoais. O representante passará o feedback do cliente para a equipe e vice-versa.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "reuniões pessoais. ", 2)
